$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for rows 2-27.
# Update it from 45330 (2024-02-08) to 45331 (2024-02-09).
for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45330) {
        $cell.Value2 = 45331
    }
}
